$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "INTU"
$ws.Range("B2").Value = 80.95999999999999
$ws.Range("C2").Value = 60.16666666666667
$ws.Range("D2").Value = 70.59999999999999

$ws.Range("A3").Value = "NVDA"
$ws.Range("B3").Value = 68.55
$ws.Range("C3").Value = 69.32380952380953
$ws.Range("D3").Value = 68.90000000000001

$ws.Range("A4").Value = "IDXX"
$ws.Range("B4").Value = 76.87
$ws.Range("C4").Value = 60.4921568627451
$ws.Range("D4").Value = 68.7

$ws.Range("A5").Value = "APH"
$ws.Range("B5").Value = 80.59
$ws.Range("C5").Value = 54.3169440546817
$ws.Range("D5").Value = 67.5

$ws.Range("A6").Value = "GE"
$ws.Range("B6").Value = 77.09
$ws.Range("C6").Value = 57.875
$ws.Range("D6").Value = 67.5

$ws.Range("A7").Value = "MSFT"
$ws.Range("B7").Value = 78.76000000000001
$ws.Range("C7").Value = 55.79749505777356
$ws.Range("D7").Value = 67.3

$ws.Range("A8").Value = "HWM"
$ws.Range("B8").Value = 73.84
$ws.Range("C8").Value = 60.59305555555556
$ws.Range("D8").Value = 67.2

$ws.Range("A9").Value = "ULTA"
$ws.Range("B9").Value = 78.56
$ws.Range("C9").Value = 53.90141260211903
$ws.Range("D9").Value = 66.2

$ws.Range("A10").Value = "NEM"
$ws.Range("B10").Value = 71.03
$ws.Range("C10").Value = 59.38267523117399
$ws.Range("D10").Value = 65.2

$ws.Range("A11").Value = "KLAC"
$ws.Range("B11").Value = 75.95999999999999
$ws.Range("C11").Value = 53.88846515318296
$ws.Range("D11").Value = 64.90000000000001

$ws.Range("A12").Value = "NFLX"
$ws.Range("B12").Value = 69.27
$ws.Range("C12").Value = 60.53185823930151
$ws.Range("D12").Value = 64.90000000000001

$ws.Range("A13").Value = "META"
$ws.Range("B13").Value = 65.7
$ws.Range("C13").Value = 62.4270698920802
$ws.Range("D13").Value = 64.09999999999999

$ws.Range("A14").Value = "MA"
$ws.Range("B14").Value = 66.84999999999999
$ws.Range("C14").Value = 59.27271133635813
$ws.Range("D14").Value = 63.1

$ws.Range("A15").Value = "RMD"
$ws.Range("B15").Value = 64.81999999999999
$ws.Range("C15").Value = 61.24901960784313
$ws.Range("D15").Value = 63

$ws.Range("A16").Value = "TEL"
$ws.Range("B16").Value = 70.94
$ws.Range("C16").Value = 52.67881280387321
$ws.Range("D16").Value = 61.8

$ws.Range("A17").Value = "CTAS"
$ws.Range("B17").Value = 65.15000000000001
$ws.Range("C17").Value = 58.47373786712579
$ws.Range("D17").Value = 61.8

$ws.Range("A18").Value = "PLTR"
$ws.Range("B18").Value = 69.45
$ws.Range("C18").Value = 53.7025641025641
$ws.Range("D18").Value = 61.6

$ws.Range("A19").Value = "V"
$ws.Range("B19").Value = 65.38
$ws.Range("C19").Value = 56.66074070361164
$ws.Range("D19").Value = 61

$ws.Range("A20").Value = "DASH"
$ws.Range("B20").Value = 66.19
$ws.Range("C20").Value = 54.57242998809394
$ws.Range("D20").Value = 60.4

$ws.Range("A21").Value = "MNST"
$ws.Range("B21").Value = 58.95
$ws.Range("C21").Value = 59.32245180363024
$ws.Range("D21").Value = 59.1

$ws.Range("A22").Value = "TT"
$ws.Range("B22").Value = 62.36
$ws.Range("C22").Value = 54.45433477188041
$ws.Range("D22").Value = 58.4

$ws.Range("A23").Value = "FAST"
$ws.Range("B23").Value = 61.38
$ws.Range("C23").Value = 54.60710172024616
$ws.Range("D23").Value = 58

$ws.Range("A24").Value = "ROL"
$ws.Range("B24").Value = 61.5
$ws.Range("C24").Value = 53.88608174588665
$ws.Range("D24").Value = 57.7

$ws.Range("A25").Value = "VRSK"
$ws.Range("B25").Value = 59.32
$ws.Range("C25").Value = 52.64907811522774
$ws.Range("D25").Value = 56

$ws.Range("A26:D26").ClearContents()
